# #324 - implement bit vectors
#
# 1) Bump the "mu name space, version 0.1.66" run pair into a single
#    "0.1.67" run (version bump).
# 2) Remove a stray empty paragraph that sits between the "`form /
#    quoted form" entry and the "`form / backquoted form" entry.
# 3) Add a new reference-card entry "#*...	bit vector" (reader syntax
#    for bit vectors), right before the existing "#x... hexadecimal
#    fixnum" entry.

$d = $word.ActiveDocument

# --- 1) version bump: 0.1.66 -> 0.1.67 -------------------------------
$d.Content.Find.Execute("0.1.66", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "0.1.67", 2) | Out-Null

# --- 2) delete the stray empty paragraph after "quoted form" ---------
# (the "`form" entry, not the later "backquoted form" / "eval backquoted
# form" / "eval-splice backquoted form" entries, whose text also ends in
# "quoted form" as a substring)
$paras = $d.Paragraphs
$quoteChar = [char]8216
$quotedFormTarget = $quoteChar + "form" + [char]9 + [char]9 + "quoted form"
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.Trim() -eq $quotedFormTarget) {
        $next = $paras.Item($i + 1)
        if ($next.Range.Text.Trim().Length -eq 0) {
            $next.Range.Delete()
        }
        break
    }
}

# --- 3) insert the new "bit vector" reference-card entry -------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*hexadecimal fixnum*") {
        $target = $p
        break
    }
}

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body>' + `
'<w:p>' + `
  '<w:pPr>' + `
    '<w:pStyle w:val="Normal"/>' + `
    '<w:ind w:hanging="0" w:start="0"/>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>' + `
      '<w:sz w:val="20"/>' + `
      '<w:szCs w:val="20"/>' + `
    '</w:rPr>' + `
  '</w:pPr>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>' + `
      '<w:i w:val="false"/>' + `
      '<w:iCs w:val="false"/>' + `
      '<w:caps w:val="false"/>' + `
      '<w:smallCaps w:val="false"/>' + `
      '<w:color w:val="24292E"/>' + `
      '<w:spacing w:val="0"/>' + `
      '<w:sz w:val="20"/>' + `
      '<w:szCs w:val="20"/>' + `
    '</w:rPr>' + `
    '<w:t xml:space="preserve">  </w:t>' + `
  '</w:r>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>' + `
      '<w:i w:val="false"/>' + `
      '<w:iCs w:val="false"/>' + `
      '<w:caps w:val="false"/>' + `
      '<w:smallCaps w:val="false"/>' + `
      '<w:color w:val="24292E"/>' + `
      '<w:spacing w:val="0"/>' + `
      '<w:sz w:val="20"/>' + `
      '<w:szCs w:val="20"/>' + `
    '</w:rPr>' + `
    '<w:t>#</w:t>' + `
  '</w:r>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>' + `
      '<w:i w:val="false"/>' + `
      '<w:iCs w:val="false"/>' + `
      '<w:caps w:val="false"/>' + `
      '<w:smallCaps w:val="false"/>' + `
      '<w:color w:val="24292E"/>' + `
      '<w:spacing w:val="0"/>' + `
      '<w:sz w:val="20"/>' + `
      '<w:szCs w:val="20"/>' + `
    '</w:rPr>' + `
    '<w:t>*</w:t>' + `
  '</w:r>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>' + `
      '<w:i w:val="false"/>' + `
      '<w:iCs w:val="false"/>' + `
      '<w:caps w:val="false"/>' + `
      '<w:smallCaps w:val="false"/>' + `
      '<w:color w:val="24292E"/>' + `
      '<w:spacing w:val="0"/>' + `
      '<w:sz w:val="20"/>' + `
      '<w:szCs w:val="20"/>' + `
    '</w:rPr>' + `
    '<w:t>...</w:t>' + `
    '<w:tab/>' + `
    '<w:tab/>' + `
  '</w:r>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/>' + `
      '<w:i w:val="false"/>' + `
      '<w:iCs w:val="false"/>' + `
      '<w:caps w:val="false"/>' + `
      '<w:smallCaps w:val="false"/>' + `
      '<w:color w:val="24292E"/>' + `
      '<w:spacing w:val="0"/>' + `
      '<w:sz w:val="16"/>' + `
      '<w:szCs w:val="16"/>' + `
    '</w:rPr>' + `
    '<w:t>bit vector</w:t>' + `
  '</w:r>' + `
'</w:p>' + `
'<w:p><w:r><w:t></w:t></w:r></w:p>' + `
'</w:body>' + `
'</w:document>' + `
'</pkg:xmlData>' + `
'</pkg:part>' + `
'</pkg:package>'

$insertPoint = $d.Range($target.Range.Start, $target.Range.Start)
$insertPoint.InsertXML($xml)
